# feat: add 2022-Q3 data
#
# Before:  Sheets = [ "总计", "2022-Q2" ]
# After:   Sheets = [ "总计", "2022-Q3" (new data), "2022-Q2" (unchanged data, moved to 3rd slot) ]
#
# - "总计" row 2 (previously the 2022-Q2 summary row) becomes the 2022-Q3
#   summary row, and a brand-new row 3 is appended holding the old 2022-Q2
#   summary values.
# - A new worksheet named "2022-Q3" is inserted between "总计" and "2022-Q2"
#   and populated with the quarter's fund-holding detail table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Row 2 used to describe 2022-Q2 (2 holdings, 0.18). Push that down to a new
# row 3 first, then overwrite row 2 with the 2022-Q3 figures.
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q2"
$totals.Range("C3").Value = 2
$totals.Range("D3").Value = 0.18

$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 2
$totals.Range("D2").Value = 0.21

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right before the existing
#    "2022-Q2" sheet, so the tab order becomes 总计, 2022-Q3, 2022-Q2.
# ---------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($oldQ2)
$q3.Name = "2022-Q3"

# Headers (row 1)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Fund rows (2 - 3). Columns D:G hold numeric-looking figures that are
# nonetheless stored as text in the source data, so force a text format
# before writing them in order to avoid silent numeric coercion.
$q3.Range("D2:G3").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "014277"
$q3.Range("C2").Value = "万家北交所慧选两年定期开放混合A"
$q3.Range("D2").Value = "3.56"
$q3.Range("E2").Value = "93.97"
$q3.Range("F2").Value = "5.07"
$q3.Range("G2").Value = "0.1805"
$q3.Range("H2").Value = 8

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "014278"
$q3.Range("C3").Value = "万家北交所慧选两年定期开放混合C"
$q3.Range("D3").Value = "0.49"
$q3.Range("E3").Value = "93.97"
$q3.Range("F3").Value = "5.07"
$q3.Range("G3").Value = "0.0248"
$q3.Range("H3").Value = 8
